$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture existing hyperlink cell info (row, address, and formatting) before
# the column insert shifts the underlying cells. The Hyperlinks collection
# itself does not track the shift, and re-Add()-ing a hyperlink stamps the
# cell with a generic "Hyperlink" style, so we snapshot + restore the real
# per-cell formatting by hand.
$links = @()
foreach ($hl in $ws.Hyperlinks) {
    $rng = $hl.Range
    $info = @{
        Row        = $rng.Row
        Address    = $hl.Address
        FontName   = $rng.Font.Name
        FontSize   = $rng.Font.Size
        Underline  = $rng.Font.Underline
        FontColor  = $rng.Font.Color
        IntPattern = $rng.Interior.Pattern
        IntColor   = $rng.Interior.Color
    }
    $links += ,$info
}

# Insert a new column before column A; this shifts columns A:J -> B:K
$ws.Columns("A").Insert()

# New column A: header + sequential row index (1-15) for each data row
$ws.Range("A1").Value = "ID for study"
for ($i = 1; $i -le 15; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $i
}

# Rebuild the hyperlinks against their new column (C, was B before the insert)
# and restore the original cell formatting that Hyperlinks.Add() overwrites.
$ws.Hyperlinks.Delete()
foreach ($info in $links) {
    $cell = $ws.Cells.Item($info.Row, 3)
    $ws.Hyperlinks.Add($cell, $info.Address) | Out-Null
    $cell.Font.Name = $info.FontName
    $cell.Font.Size = $info.FontSize
    $cell.Font.Underline = $info.Underline
    $cell.Font.Color = $info.FontColor
    if ($info.IntPattern -ne -4142) {
        $cell.Interior.Color = $info.IntColor
    }
}

# Column J (originally I) was manually widened
$ws.Columns("J").ColumnWidth = 58.33

# View state: zoom to 85%, selection on C16
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("C16").Select() | Out-Null
